$d = $word.ActiveDocument

# Helper: append one or more runs of plain text followed by a colored run,
# all at the very end of a paragraph's content (i.e. just before its
# paragraph mark). $segments is an array of hashtables:
#   @{ Text = "..."; Color = "548DD4" }   -> colored run
#   @{ Text = "..." }                      -> plain run (no color)
function Append-Runs($paraIndex, $segments) {
    $p = $d.Paragraphs($paraIndex)
    $pos = $p.Range.End - 1
    foreach ($seg in $segments) {
        $text = $seg.Text
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter($text)
        $newPos = $pos + $text.Length
        if ($seg.ContainsKey("Color")) {
            $fr = $d.Range($pos, $newPos)
            $fr.Font.Color = $seg.Color
        }
        $pos = $newPos
    }
}

# Helper: within a paragraph, replace the first occurrence of $oldText with
# $newText (plain, no formatting change) using Find/Replace scoped to that
# paragraph's range.
function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range.Duplicate
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- Para 3: "Ingresar al Sistema " -> + " " + colored "1ero"
Append-Runs 3 @(
    @{ Text = " " },
    @{ Text = "1ero"; Color = "548DD4" }
)

# --- Para 4: "Registro de Cliente" -> + " " + colored "1ero"
Append-Runs 4 @(
    @{ Text = " " },
    @{ Text = "1ero"; Color = "548DD4" }
)

# --- Para 5: "Registro Pedido Cliente" -> + " " + colored "1ero"
Append-Runs 5 @(
    @{ Text = " " },
    @{ Text = "1ero"; Color = "548DD4" }
)

# --- Para 7: "Registro Materia Prima (Cantidad / Costo)   3 3 3 3  33"
#     -> trailing junk "   3 3 3 3  33" becomes "  ", then colored " " + "3ero"
Replace-InParagraph 7 "   3 3 3 3  33" "  "
Append-Runs 7 @(
    @{ Text = " "; Color = "FF0000" },
    @{ Text = "3ero"; Color = "FF0000" }
)

# --- Para 8: "Registro Mercancía 3 3  3 3 3 3 " (bookmark follows)
#     -> trailing junk " 3 3  3 3 3 3 " becomes " ", then colored "3ero"
Replace-InParagraph 8 " 3 3  3 3 3 3 " " "
Append-Runs 8 @(
    @{ Text = "3ero"; Color = "FF0000" }
)

# --- Para 9: "Registro proveedores " -> + " " + colored "1ero "
Append-Runs 9 @(
    @{ Text = " " },
    @{ Text = "1ero "; Color = "548DD4" }
)

# --- Para 11: "Registro Empleados" -> + " " + colored "1ero" + " "
Append-Runs 11 @(
    @{ Text = " " },
    @{ Text = "1ero"; Color = "548DD4" },
    @{ Text = " " }
)

# --- Para 17: "Modificación estatus Empleado" -> + "  " + colored "2do"
Append-Runs 17 @(
    @{ Text = "  " },
    @{ Text = "2do"; Color = "548DD4" }
)

# --- Para 18: "Modificación estatus Cliente" -> + " " + colored "2do"
Append-Runs 18 @(
    @{ Text = " " },
    @{ Text = "2do"; Color = "548DD4" }
)

# --- Para 19: "Modificación estatus Proveedor" -> + " " + colored "2do"
Append-Runs 19 @(
    @{ Text = " " },
    @{ Text = "2do"; Color = "548DD4" }
)

# --- Para 20: "Eliminación Empleado" -> + " " + colored "2do" + " "
Append-Runs 20 @(
    @{ Text = " " },
    @{ Text = "2do"; Color = "548DD4" },
    @{ Text = " " }
)

# --- Para 21: "Eliminación Cliente" -> + colored " 2do "
Append-Runs 21 @(
    @{ Text = " 2do "; Color = "548DD4" }
)

# --- Para 22: "Eliminación Proveedor " -> + colored "2do "
Append-Runs 22 @(
    @{ Text = "2do "; Color = "548DD4" }
)
